$wb = $excel.ActiveWorkbook

# --- Product sheet: swap product image paths from .png to .jpg ---
$wsProduct = $wb.Worksheets.Item("Product")
for ($r = 2; $r -le 30; $r += 2) {
    $wsProduct.Range("D$r").Value = "/Image/Products/1.jpg"
}
for ($r = 3; $r -le 31; $r += 2) {
    $wsProduct.Range("D$r").Value = "/Image/Products/2.jpg"
}

# --- ProductItem_VariationOptions sheet: remove stray styled cells in column G ---
$wsPIVO = $wb.Worksheets.Item("ProductItem_VariationOptions")
$wsPIVO.Range("G8:G20").Clear()

# --- Update active sheet / selection to reflect the saved view state ---
$wsProduct.Activate()
$wsProduct.Range("G8").Select()
